# The calibration data rows (A2:D18) need to be re-sorted in ascending
# order of the "time (s)" column (column A) — newly-collected samples
# were appended out of chronological order during calibration.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A2:D18")
$sortKey   = $ws.Range("A2:A18")

# xlAscending = 1, xlNo (no header in selected range) = 2
$dataRange.Sort($sortKey, 1, $null, $null, 1, $null, 1, 2)
